# Append the profit-data row for 2026-01-21 (run date 2026-01-21) as the
# new last row (row 58) of the worksheet, matching the existing layout:
#   A = Date (text, e.g. "01/21/2026")
#   B = Portfolio Value(USD)
#   C = BTC %
#   D = KAS %
#   E = KAS Profit(USD)
#   F = KAS Profit(%)
#   G = KAS Total Profit(USD)
#   H = KAS Total Profit(%)
#   I = BTC Profit(USD)
#   J = BTC Profit(%)
#   K = Combined Total Profit(USD)
#   L = Combined Total Profit(%)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 58

# Column A holds the date as plain text in this workbook (not a real Excel
# date value). Assigning a date-looking string normally makes Excel infer
# a date and reformat the cell, so force the cell to Text first, write the
# value, then strip the formatting we had to add back off again so the
# cell ends up as a plain, unstyled text cell - exactly like the rest of
# the column.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat() = "@"
$dateCell.Value() = "01/21/2026"
$dateCell.ClearFormats()

$ws.Cells.Item($row, 2).Value()  = 11599.77
$ws.Cells.Item($row, 3).Value()  = 0.236403919765626
$ws.Cells.Item($row, 4).Value()  = 0.763596080234374
$ws.Cells.Item($row, 5).Value()  = -194.95
$ws.Cells.Item($row, 6).Value()  = -28.46
$ws.Cells.Item($row, 7).Value()  = -21884.22
$ws.Cells.Item($row, 8).Value()  = -71.19
$ws.Cells.Item($row, 9).Value()  = -460.62
$ws.Cells.Item($row, 10).Value() = -14.38
$ws.Cells.Item($row, 11).Value() = -22344.84
$ws.Cells.Item($row, 12).Value() = -65.83
